$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2G")

# Set custom width for column C (target stored width 39.5546875 characters)
$ws.Columns.Item(3).ColumnWidth = 38.6666666666667

# New row 16 with the new shared string text, matching the wrap-text style
# used by similar rows (e.g. row 13), and a row height of 28.8 (two wrapped lines).
$ws.Range("A16").Value = "Are the SIM card installed ?" + [char]10 + "test indirectly with the aid of the led"
$ws.Range("A16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 28.8
